$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44357
$ws.Range("J2").Value = 340
$ws.Range("K2").Value = 28000
$ws.Range("L2").Value = 30000
$ws.Range("M2").Value = 29000
$ws.Range("P2").Value = 1160

# Row 3
$ws.Range("D3").Value = 44455
$ws.Range("J3").Value = 800

# Row 4
$ws.Range("D4").Value = 44461
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 23000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24000
$ws.Range("P4").Value = 960

# Row 5
$ws.Range("D5").Value = 44489
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = 19000
$ws.Range("P5").Value = 760

# Row 6
$ws.Range("D6").Value = 44398
$ws.Range("K6").Value = 26000
$ws.Range("L6").Value = 28000
$ws.Range("M6").Value = 27000
$ws.Range("P6").Value = 1080

# Row 7
$ws.Range("D7").Value = 44364
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 28000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 29000
$ws.Range("P7").Value = 1160

# Row 8
$ws.Range("D8").Value = 44482
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 18000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19000
$ws.Range("P8").Value = 760

# Row 9
$ws.Range("D9").Value = 44384
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 26000
$ws.Range("L9").Value = 28000
$ws.Range("M9").Value = 27000
$ws.Range("P9").Value = 1080

# Row 10
$ws.Range("D10").Value = 44391
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 26000
$ws.Range("L10").Value = 28000
$ws.Range("M10").Value = 27000
$ws.Range("P10").Value = 1080

# Row 11
$ws.Range("D11").Value = 44503
$ws.Range("K11").Value = 11000
$ws.Range("L11").Value = 13000
$ws.Range("M11").Value = 12000
$ws.Range("P11").Value = 480

# Row 12
$ws.Range("D12").Value = 44448
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 28000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 29000
$ws.Range("P12").Value = 1160

# Row 13
$ws.Range("D13").Value = 44413
$ws.Range("J13").Value = 700
$ws.Range("K13").Value = 26000
$ws.Range("L13").Value = 28000
$ws.Range("M13").Value = 27000
$ws.Range("P13").Value = 1080

# Row 15
$ws.Range("D15").Value = 44419
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 27000
$ws.Range("L15").Value = 29000
$ws.Range("M15").Value = 28000
$ws.Range("P15").Value = 1120

# Row 16
$ws.Range("D16").Value = 44490
$ws.Range("K16").Value = 16000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 17000
$ws.Range("P16").Value = 680

# Row 17
$ws.Range("D17").Value = 44497
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 13000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14000
$ws.Range("P17").Value = 560

# Row 18
$ws.Range("D18").Value = 44392
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 26000
$ws.Range("L18").Value = 28000
$ws.Range("M18").Value = 27000
$ws.Range("P18").Value = 1080

# Row 19
$ws.Range("D19").Value = 44433
$ws.Range("J19").Value = 400

# Row 20
$ws.Range("D20").Value = 44377
$ws.Range("J20").Value = 500

# Row 21
$ws.Range("D21").Value = 44350
$ws.Range("J21").Value = 700
$ws.Range("K21").Value = 28000
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = 29000
$ws.Range("P21").Value = 1160

# Row 22
$ws.Range("D22").Value = 44476
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 23000
$ws.Range("L22").Value = 24000
$ws.Range("M22").Value = 23500
$ws.Range("P22").Value = 940

# Row 23
$ws.Range("D23").Value = 44399
$ws.Range("K23").Value = 26000
$ws.Range("L23").Value = 28000
$ws.Range("M23").Value = 27000
$ws.Range("P23").Value = 1080

# Row 25
$ws.Range("D25").Value = 44469
$ws.Range("J25").Value = 600
$ws.Range("K25").Value = 22000
$ws.Range("L25").Value = 24000
$ws.Range("M25").Value = 23000
$ws.Range("P25").Value = 920

# Row 26
$ws.Range("D26").Value = 44475
$ws.Range("J26").Value = 1000
$ws.Range("K26").Value = 22000
$ws.Range("L26").Value = 24000
$ws.Range("M26").Value = 23000
$ws.Range("P26").Value = 920

# Row 27
$ws.Range("D27").Value = 44356
$ws.Range("J27").Value = 300
$ws.Range("K27").Value = 26000
$ws.Range("M27").Value = 27000
$ws.Range("P27").Value = 1080

# Row 28
$ws.Range("D28").Value = 44363
$ws.Range("J28").Value = 240
$ws.Range("K28").Value = 28000
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = 29000
$ws.Range("P28").Value = 1160

# Row 29
$ws.Range("D29").Value = 44483
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 18000
$ws.Range("L29").Value = 20000
$ws.Range("M29").Value = 19000
$ws.Range("P29").Value = 760

# Row 30
$ws.Range("D30").Value = 44462
$ws.Range("J30").Value = 400
$ws.Range("K30").Value = 22000
$ws.Range("L30").Value = 23000
$ws.Range("M30").Value = 22500
$ws.Range("P30").Value = 900

# Row 31
$ws.Range("D31").Value = 44405
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 26000
$ws.Range("L31").Value = 28000
$ws.Range("M31").Value = 27000
$ws.Range("P31").Value = 1080

# Row 32
$ws.Range("D32").Value = 44412
$ws.Range("J32").Value = 600
$ws.Range("K32").Value = 25000
$ws.Range("L32").Value = 27000
$ws.Range("M32").Value = 26000
$ws.Range("P32").Value = 1040

# Row 33
$ws.Range("D33").Value = 44468
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 23000
$ws.Range("L33").Value = 25000
$ws.Range("M33").Value = 24000
$ws.Range("P33").Value = 960

# Row 34
$ws.Range("D34").Value = 44370
$ws.Range("K34").Value = 27000
$ws.Range("M34").Value = 27500
$ws.Range("P34").Value = 1100

# Row 35
$ws.Range("D35").Value = 44504
$ws.Range("J35").Value = 600
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 13000
$ws.Range("M35").Value = 12000
$ws.Range("P35").Value = 480

# Row 36
$ws.Range("D36").Value = 44447
$ws.Range("J36").Value = 600
$ws.Range("K36").Value = 28000
$ws.Range("L36").Value = 30000
$ws.Range("M36").Value = 29000
$ws.Range("P36").Value = 1160

# Row 37
$ws.Range("D37").Value = 44434
$ws.Range("J37").Value = 500
$ws.Range("K37").Value = 28000
$ws.Range("L37").Value = 30000
$ws.Range("M37").Value = 29000
$ws.Range("P37").Value = 1160

# Row 38
$ws.Range("D38").Value = 44385
$ws.Range("J38").Value = 500

# Row 39
$ws.Range("D39").Value = 44371
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 28000
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = 29000
$ws.Range("P39").Value = 1160

# Row 40
$ws.Range("D40").Value = 44406
$ws.Range("K40").Value = 26000
$ws.Range("L40").Value = 28000
$ws.Range("M40").Value = 27000
$ws.Range("P40").Value = 1080

# Row 41
$ws.Range("D41").Value = 44427
$ws.Range("J41").Value = 300
$ws.Range("K41").Value = 28000
$ws.Range("L41").Value = 30000
$ws.Range("M41").Value = 29000
$ws.Range("P41").Value = 1160

# Row 42
$ws.Range("D42").Value = 44441
$ws.Range("J42").Value = 700
$ws.Range("K42").Value = 28000
$ws.Range("L42").Value = 30000
$ws.Range("M42").Value = 29000
$ws.Range("P42").Value = 1160

# Row 43
$ws.Range("D43").Value = 44426
$ws.Range("J43").Value = 400

# Row 44
$ws.Range("D44").Value = 44343
$ws.Range("J44").Value = 200

# Row 45
$ws.Range("D45").Value = 44435
$ws.Range("J45").Value = 900
$ws.Range("K45").Value = 28000
$ws.Range("L45").Value = 30000
$ws.Range("M45").Value = 29000
$ws.Range("P45").Value = 1160

# Row 46
$ws.Range("D46").Value = 44349
$ws.Range("J46").Value = 600
$ws.Range("K46").Value = 26000
$ws.Range("L46").Value = 28000
$ws.Range("M46").Value = 27000
$ws.Range("P46").Value = 1080
